$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.8632254154274
$ws.Range("C2").Value = 11.71562406796107
$ws.Range("E2").Value = 16.88187714018735
$ws.Range("F2").Value = 34.2423460306265
$ws.Range("G2").Value = 24.7555402894241
$ws.Range("H2").Value = 13.26614312731674
$ws.Range("I2").Value = 20.93096524664174
$ws.Range("J2").Value = 7.031460858160722
$ws.Range("L2").Value = 13.12542293019256
$ws.Range("N2").Value = 16.99889466005951
$ws.Range("O2").Value = 19.65991991917291
$ws.Range("B3").Value = 14.35589679774867
$ws.Range("C3").Value = 11.63343596916445
$ws.Range("E3").Value = 16.91007058029157
$ws.Range("F3").Value = 34.24985981617488
$ws.Range("G3").Value = 24.76097453403908
$ws.Range("H3").Value = 13.30705496182806
$ws.Range("I3").Value = 21.02734140717243
$ws.Range("J3").Value = 7.021348884147555
$ws.Range("L3").Value = 13.09343676531269
$ws.Range("N3").Value = 17.02924377958301
$ws.Range("O3").Value = 19.71659045715352
$ws.Range("B4").Value = 14.03612098167844
$ws.Range("C4").Value = 11.5826016541428
$ws.Range("E4").Value = 16.92981617082593
$ws.Range("F4").Value = 34.26283741874862
$ws.Range("G4").Value = 24.77326584739384
$ws.Range("H4").Value = 13.33437981971176
$ws.Range("I4").Value = 21.0906699921699
$ws.Range("J4").Value = 7.015104556828565
$ws.Range("L4").Value = 13.07554667894149
$ws.Range("N4").Value = 17.04955219004174
$ws.Range("O4").Value = 19.75592606303271
$ws.Range("B5").Value = 13.90391811629674
$ws.Range("C5").Value = 11.56180373249875
$ws.Range("E5").Value = 16.93847517380434
$ws.Range("F5").Value = 34.27022965163775
$ws.Range("G5").Value = 24.78052206658044
$ws.Range("H5").Value = 13.34606910890986
$ws.Range("I5").Value = 21.11752064457685
$ws.Range("J5").Value = 7.012551348493265
$ws.Range("L5").Value = 13.06870127360822
$ws.Range("N5").Value = 17.05824964274589
$ws.Range("O5").Value = 19.77309453897431
$ws.Range("B6").Value = 13.88185754715803
$ws.Range("C6").Value = 11.55834557929803
$ws.Range("E6").Value = 16.93994999616776
$ws.Range("F6").Value = 34.2715841799597
$ws.Range("G6").Value = 24.78186252962882
$ws.Range("H6").Value = 13.34804357105051
$ws.Range("I6").Value = 21.12204219672368
$ws.Range("J6").Value = 7.012126880984401
$ws.Range("L6").Value = 13.06759161306716
$ws.Range("N6").Value = 17.05971933260892
$ws.Range("O6").Value = 19.77601405747949
$ws.Range("B7").Value = 14.03434544419671
$ws.Range("C7").Value = 11.58232148621528
$ws.Range("E7").Value = 16.92993046878184
$ws.Range("F7").Value = 34.26292859515328
$ws.Range("G7").Value = 24.77335461480095
$ws.Range("H7").Value = 13.33453522198823
$ws.Range("I7").Value = 21.091027884198
$ws.Range("J7").Value = 7.015070157720653
$ws.Range("L7").Value = 13.07545255145019
$ws.Range("N7").Value = 17.04966777893348
$ws.Range("O7").Value = 19.75615299537125
$ws.Range("B8").Value = 14.69012372187811
$ws.Range("C8").Value = 11.68736888912771
$ws.Range("E8").Value = 16.89109310663236
$ws.Range("F8").Value = 34.24320148084546
$ws.Range("G8").Value = 24.75555347580696
$ws.Range("H8").Value = 13.27979177973602
$ws.Range("I8").Value = 20.96333336395905
$ws.Range("J8").Value = 7.027981828946603
$ws.Range("L8").Value = 13.11403408430327
$ws.Range("N8").Value = 17.00901202475047
$ws.Range("O8").Value = 19.67851617592697
$ws.Range("B9").Value = 15.90303355660076
$ws.Range("C9").Value = 11.88996198846041
$ws.Range("E9").Value = 16.83423722324098
$ws.Range("F9").Value = 34.27081012778373
$ws.Range("G9").Value = 24.79179814990458
$ws.Range("H9").Value = 13.18994452025986
$ws.Range("I9").Value = 20.7459068516285
$ws.Range("J9").Value = 7.053005273709004
$ws.Range("L9").Value = 13.20334110665529
$ws.Range("N9").Value = 16.94253811110016
$ws.Range("O9").Value = 19.56240240786058
$ws.Range("B10").Value = 16.74086890661423
$ws.Range("C10").Value = 12.03611375743128
$ws.Range("E10").Value = 16.80421194250793
$ws.Range("F10").Value = 34.33133957647747
$ws.Range("G10").Value = 24.86180794879157
$ws.Range("H10").Value = 13.1346183527667
$ws.Range("I10").Value = 20.60630892491751
$ws.Range("J10").Value = 7.07119103262261
$ws.Range("L10").Value = 13.27694483931374
$ws.Range("N10").Value = 16.90173845888887
$ws.Range("O10").Value = 19.49926332470065
$ws.Range("B11").Value = 17.10892248817138
$ws.Range("C11").Value = 12.10188680349506
$ws.Range("E11").Value = 16.79309735018839
$ws.Range("F11").Value = 34.36755252995965
$ws.Range("G11").Value = 24.90303318858504
$ws.Range("H11").Value = 13.11177180570197
$ws.Range("I11").Value = 20.54718570010252
$ws.Range("J11").Value = 7.079417036458264
$ws.Range("L11").Value = 13.31208717734023
$ws.Range("N11").Value = 16.88491477522829
$ws.Range("O11").Value = 19.47538023510961
$ws.Range("B12").Value = 17.24630586149665
$ws.Range("C12").Value = 12.1266799918149
$ws.Range("E12").Value = 16.7892537228853
$ws.Range("F12").Value = 34.38250579154769
$ws.Range("G12").Value = 24.91998535670103
$ws.Range("H12").Value = 13.10345448531418
$ws.Range("I12").Value = 20.52542813703625
$ws.Range("J12").Value = 7.082524904699997
$ws.Range("L12").Value = 13.32562631971414
$ws.Range("N12").Value = 16.87879308653035
$ws.Range("O12").Value = 19.46703402023193
$ws.Range("B13").Value = 17.2168080235244
$ws.Range("C13").Value = 12.12134556664179
$ws.Range("E13").Value = 16.7900652837503
$ws.Range("F13").Value = 34.37923031672505
$ws.Range("G13").Value = 24.91627490084341
$ws.Range("H13").Value = 13.10523090169459
$ws.Range("I13").Value = 20.53008592884207
$ws.Range("J13").Value = 7.081855893897811
$ws.Range("L13").Value = 13.32270024383878
$ws.Range("N13").Value = 16.88010043358613
$ws.Range("O13").Value = 19.46880046433687
$ws.Range("B14").Value = 17.12026557327811
$ws.Range("C14").Value = 12.10392889536424
$ws.Range("E14").Value = 16.79277381672167
$ws.Range("F14").Value = 34.36875793683145
$ws.Range("G14").Value = 24.90440102688053
$ws.Range("H14").Value = 13.11108083564687
$ws.Range("I14").Value = 20.5453830367236
$ws.Range("J14").Value = 7.079672868149302
$ws.Range("L14").Value = 13.31319645303008
$ws.Range("N14").Value = 16.88440615181232
$ws.Range("O14").Value = 19.47467959201469
$ws.Range("B15").Value = 17.06086829128389
$ws.Range("C15").Value = 12.09324555119307
$ws.Range("E15").Value = 16.7944804153712
$ws.Range("F15").Value = 34.3625045629529
$ws.Range("G15").Value = 24.89730232548283
$ws.Range("H15").Value = 13.11470761747671
$ws.Range("I15").Value = 20.55483518356859
$ws.Range("J15").Value = 7.078334759534859
$ws.Range("L15").Value = 13.30740503815789
$ws.Range("N15").Value = 16.88707594769339
$ws.Range("O15").Value = 19.47837165370209
$ws.Range("B16").Value = 16.71654284072833
$ws.Range("C16").Value = 12.03180004129929
$ws.Range("E16").Value = 16.80498945045962
$ws.Range("F16").Value = 34.32914697014788
$ws.Range("G16").Value = 24.8593018979018
$ws.Range("H16").Value = 13.13615817332106
$ws.Range("I16").Value = 20.61026101651017
$ws.Range("J16").Value = 7.070652467279798
$ws.Range("L16").Value = 13.27468099593243
$ws.Range("N16").Value = 16.90287281330654
$ws.Range("O16").Value = 19.50092168747782
$ws.Range("B17").Value = 16.50187759990274
$ws.Range("C17").Value = 11.99391478794676
$ws.Range("E17").Value = 16.81208759419215
$ws.Range("F17").Value = 34.31090073401174
$ws.Range("G17").Value = 24.83838708410535
$ws.Range("H17").Value = 13.14991225124001
$ws.Range("I17").Value = 20.64538577905035
$ws.Range("J17").Value = 7.06592735442437
$ws.Range("L17").Value = 13.2550259280856
$ws.Range("N17").Value = 16.91300795505777
$ws.Range("O17").Value = 19.51599627925019
$ws.Range("B18").Value = 16.37718213886475
$ws.Range("C18").Value = 11.97205782864069
$ws.Range("E18").Value = 16.81640975911917
$ws.Range("F18").Value = 34.30122331979911
$ws.Range("G18").Value = 24.82724066732085
$ws.Range("H18").Value = 13.15804173487325
$ws.Range("I18").Value = 20.66600085941816
$ws.Range("J18").Value = 7.063205215561648
$ws.Range("L18").Value = 13.24387751133865
$ws.Range("N18").Value = 16.91900088348367
$ws.Range("O18").Value = 19.52512217341728
$ws.Range("B19").Value = 16.33475543898701
$ws.Range("C19").Value = 11.96464638401679
$ws.Range("E19").Value = 16.81791432135184
$ws.Range("F19").Value = 34.29808730387565
$ws.Range("G19").Value = 24.82361856740466
$ws.Range("H19").Value = 13.16083175600218
$ws.Range("I19").Value = 20.67305153922712
$ws.Range("J19").Value = 7.062282806500657
$ws.Range("L19").Value = 13.24012996854074
$ws.Range("N19").Value = 16.9210580793028
$ws.Range("O19").Value = 19.52829019524343
$ws.Range("B20").Value = 16.52485677062522
$ws.Range("C20").Value = 11.9979546879903
$ws.Range("E20").Value = 16.81130720139736
$ws.Range("F20").Value = 34.31275853681814
$ws.Range("G20").Value = 24.84052213727744
$ws.Range("H20").Value = 13.1484254896195
$ws.Range("I20").Value = 20.64160401364686
$ws.Range("J20").Value = 7.0664308069246
$ws.Range("L20").Value = 13.25710208097624
$ws.Range("N20").Value = 16.91191213846146
$ws.Range("O20").Value = 19.51434441762294
$ws.Range("B21").Value = 17.1486772591326
$ws.Range("C21").Value = 12.10904776687224
$ws.Range("E21").Value = 16.79196834800742
$ws.Range("F21").Value = 34.37180033748269
$ws.Range("G21").Value = 24.90785234519025
$ws.Range("H21").Value = 13.10935349640768
$ws.Range("I21").Value = 20.54087276719443
$ws.Range("J21").Value = 7.080314272850882
$ws.Range("L21").Value = 13.31598172222917
$ws.Range("N21").Value = 16.88313470294758
$ws.Range("O21").Value = 19.47293379766402
$ws.Range("B22").Value = 17.54473611206609
$ws.Range("C22").Value = 12.18098603769272
$ws.Range("E22").Value = 16.7814577736709
$ws.Range("F22").Value = 34.41761193435774
$ws.Range("G22").Value = 24.95966830894307
$ws.Range("H22").Value = 13.08576563490836
$ws.Range("I22").Value = 20.47871830363671
$ws.Range("J22").Value = 7.089346101984835
$ws.Range("L22").Value = 13.35580883338824
$ws.Range("N22").Value = 16.86577850802484
$ws.Range("O22").Value = 19.44993755649949
$ws.Range("B23").Value = 17.33445020648726
$ws.Range("C23").Value = 12.14265598029994
$ws.Range("E23").Value = 16.78687293065187
$ws.Range("F23").Value = 34.39250326213829
$ws.Range("G23").Value = 24.93130138339523
$ws.Range("H23").Value = 13.09817657978115
$ws.Range("I23").Value = 20.51155425462737
$ws.Range("J23").Value = 7.08452960915401
$ws.Range("L23").Value = 13.33443164066229
$ws.Range("N23").Value = 16.87490921958342
$ws.Range("O23").Value = 19.46183828828052
$ws.Range("B24").Value = 16.51447187084181
$ws.Range("C24").Value = 11.99612848575285
$ws.Range("E24").Value = 16.8116592651686
$ws.Range("F24").Value = 34.31191609246273
$ws.Range("G24").Value = 24.83955414496342
$ws.Range("H24").Value = 13.14909696280016
$ws.Range("I24").Value = 20.64331243727932
$ws.Range("J24").Value = 7.066203213487185
$ws.Range("L24").Value = 13.25616297967623
$ws.Range("N24").Value = 16.91240703970587
$ws.Range("O24").Value = 19.51508979368768
$ws.Range("B25").Value = 15.5836983085088
$ws.Range("C25").Value = 11.83558357157838
$ws.Range("E25").Value = 16.84755353376477
$ws.Range("F25").Value = 34.25625643092198
$ws.Range("G25").Value = 24.77436674883826
$ws.Range("H25").Value = 13.21237517643256
$ws.Range("I25").Value = 20.80119169987469
$ws.Range("J25").Value = 7.046269356000654
$ws.Range("L25").Value = 13.17775287118515
$ws.Range("N25").Value = 16.94253811110016
$ws.Range("O25").Value = 19.58993056838318
